$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update J4 to the new title value (Mumbai instead of Bangalore)
$ws.Range("J4").Value = "title=Used cars in Mumbai - GoZoomo"

# Add new row 5 data
$ws.Range("A5").Value = "ts4"
$ws.Range("D5").Value = "com.selenium.test"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = "get_alert"
$ws.Range("J5").Value = "test=Login or register to get alerts."

# Update the selection to match the diff (activeCell A2, sqref row 2)
$ws.Range("A2").Select()
$ws.Rows.Item(2).Select()
